$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 (shifts rows 6-14 down to 7-15)
$ws.Rows(6).Insert()

# Append two new TODO rows at the bottom first, so the new shared-string
# entries are created in the same order as the target workbook
# (Auto-orient, then "must have 4 modes", then "Add credits").
$ws.Range("A17").Value = "Auto-orient by EXIF + reset exif orient tag (make sure other metadata is left intact)"
$ws.Range("A17").Style = "Good"

# New TODO item in the freshly inserted row (column B, matching the adjacent
# "nicer source sync mode selection" item above it in B5)
$ws.Range("B6").Value = "must have 4 modes = one is ""do not modify - just renaming and adding to output mix"""

$ws.Range("A16").Value = "Add credits for CCR-Exif and NativeJpg"

# Mark "Tools doesn't need global settings..." (now in A8) as Explanatory Text
$ws.Range("A8").Style = "Explanatory Text"

# Mark "Make Source name edit focused control..." (now in A12) and
# "Add seconds to default pattern..." (now in A14) as Good (done items)
$ws.Range("A12").Style = "Good"
$ws.Range("A14").Style = "Good"

# Update selection to match the new last-edited cell
$ws.Range("A16").Select() | Out-Null
